# resolve_report_task.xlsx uplift (pyxform/cht-conf update)
#
# Summary of the edit:
#  - "survey" sheet: drop the now-unused label column values ("NO_LABEL") for
#    every question row, and rename the question type "text" -> "hidden" for
#    those same rows (begin_group/end_group rows keep their type untouched).
#  - "survey" sheet: tidy up the conditional-formatting ranges so they cover
#    the real data range (A2:G10006 / C2:C10006) instead of the old
#    "33/34"-row split, and extend the label-required rule so it also
#    tolerates type="hidden" rows and rows that already have a calculation.
#  - "settings" sheet: remove the obsolete form_id column entirely (shifting
#    version/style/namespaces one column to the left) and update the
#    corresponding cell comments to match their new column.
#  - Selection bookkeeping to match the authors' final cursor position.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# ---------------------------------------------------------------------------
# 1. "survey" sheet: remove the NO_LABEL placeholder labels (column C) and
#    rename type "text" -> "hidden" on the rows that had them.
# ---------------------------------------------------------------------------

$hiddenRows = @(4, 5, 7, 8, 10, 12, 13)
foreach ($r in $hiddenRows) {
    $survey.Cells.Item($r, 1).Value = "hidden"
}

$labelRows = @(3, 4, 5, 6, 7, 8, 10, 11, 12, 13)
foreach ($r in $labelRows) {
    $survey.Cells.Item($r, 3).ClearContents()
}

# ---------------------------------------------------------------------------
# 2. "survey" sheet: simplify / re-scope the conditional formatting.
# ---------------------------------------------------------------------------

$cf = $survey.Cells.FormatConditions

# The last six rules only ever applied to the stray "C33" cell - drop them.
for ($i = $cf.Count; $i -ge 9; $i--) {
    $cf.Item($i).Delete()
}

# Rules 1-5 (begin_group / end_group / begin_repeat / end_repeat / blank-A
# highlighting) now cover the whole data range.
for ($i = 1; $i -le 5; $i++) {
    $cf.Item($i).ModifyAppliesToRange($survey.Range("A2:G10006"))
}

# Rule 7 (required-label-on-column-C) now also covers the whole data range,
# and is relaxed to allow type="hidden" rows or rows that already specify a
# calculation (column F) instead of a label.
$cf.Item(7).ModifyAppliesToRange($survey.Range("C2:C10006"))
$cf.Item(7).Formula1 = 'AND(ISBLANK(C2),NOT(OR(ISBLANK($A2),$A2="calculate",$A2="hidden")),ISBLANK($F2))'

# ---------------------------------------------------------------------------
# 3. "settings" sheet: drop the form_id column (B), shifting version/style/
#    namespaces one column to the left, and fix up the comments that
#    describe each column so they still match the right cell.
# ---------------------------------------------------------------------------

# Re-point the comments that will survive the shift to their new content
# *before* deleting the column (comments don't follow a Range.Delete shift).
$settings.Range("C1").Comment.Text('Set to ‘pages’ to indicate that groups with the `field-list` appearance represent separate form pages (and all other questions will be shown on their own page). ')
$settings.Range("D1").Comment.Text('Custom namespaces supported in the form.  `cht` must be included here to use the custom `instance::cht` columns on the survey sheet.')
$settings.Range("E1").Comment.Delete()

$settings.Range("B1:B2").Delete()

$settings.Range("B1").Comment.Text('The unique version code that identifies the current state of the form. A common convention is to use a format like yyyymmddrr. For example, 2017021501 is the 1st revision from Feb 15th, 2017.

By default, this template uses a formula to create a date-based version that will update automatically.')

# ---------------------------------------------------------------------------
# 4. Restore the selections the authors left behind.
# ---------------------------------------------------------------------------

$settings.Activate()
$settings.Range("B1").Select()

$survey.Activate()
$survey.Range("C9").Select()
